$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New C:D values resulting from the Kbar now-casting bug fix
$values = @(
    @(16905.43710002151, 12642.0253427297),
    @(15167.54131043129, 17422.98034507523),
    @(22274.84716920583, 24175.56059344363),
    @(24351.20418816822, 27392.09333640726),
    @(46186.82982560474, 35169.86263232183),
    @(48541.71257420899, 42426.86445458703),
    @(33661.73353113382, 31338.88368232083),
    @(39658.43649622773, 44305.51895392578),
    @(36589.3786654189, 38371.14186510815),
    @(31299.90300849402, 59476.95151094632),
    @(55486.3589472441, 49745.67670599426),
    @(49923.77865149324, 56070.05567294502),
    @(94234.57181599503, 80223.16299219802),
    @(67448.82948719188, 63576.41873186399),
    @(47550.64933182385, 73343.09185623107),
    @(43526.59953065449, 80387.06083646334),
    @(96455.76432936463, 71433.45596493836),
    @(144035.951052152, 83327.35382016608),
    @(82885.33981406162, 72900.50691780023),
    @(68894.66952040972, 76739.36906774611),
    @(69537.9363591833, 77850.66476853276),
    @(68112.40602214848, 86107.25030310568),
    @(110501.1518150782, 85876.91262145471),
    @(87370.715998074, 82852.21424654043),
    @(87071.75839233374, 85738.67527702558),
    @(96537.67647015587, 101479.0392176395),
    @(88339.75848957786, 92542.00761559447),
    @(75778.562372771, 98098.28184531824),
    @(117871.5604828493, 107025.7918995116),
    @(92765.75289327549, 99567.00452558097),
    @(79996.31761524217, 86009.57063899832),
    @(90819.52290872802, 106840.8799625777),
    @(99166.04664380381, 120926.6294421199),
    @(100365.0290816926, 121407.0812581186),
    @(90059.39297763631, 108694.5912522921),
    @(161915.8497850943, 129907.8190418039),
    @(255293.6133246898, 128008.4324912184),
    @(141424.7649926382, 135874.7051100314),
    @(157970.286754904, 128921.5677877982),
    @(128130.1754669283, 133007.2344029587),
    @(133655.1957826103, 138102.3704210776),
    @(129085.50828987, 136262.517368438),
    @(157226.6537002815, 138340.8434067073),
    @(131884.6424570323, 135722.3289850291),
    @(192460.5974124373, 179695.0437171199),
    @(184787.2901212575, 156504.7498849096),
    @(129857.7584470511, 162296.2522789774),
    @(112696.920647212, 142539.5738979665),
    @(241918.2639353328, 240077.3497923528)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 3).Value = $values[$i][0]
    $ws.Cells.Item($row, 4).Value = $values[$i][1]
}

$wb.Save()